$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing values: A2 "Profile" -> "Sheet", A3 "Sheet" -> "Tube"
$ws.Range("A2").Value = "Sheet"
$ws.Range("A3").Value = "Tube"

# Remove the old A4 ("Pipe") row entirely
$ws.Range("A4").ClearContents()

# Update the selection to reflect the new last cell
$ws.Range("A4").Select()
